# Lesson05 - Part IV - Finished Excel Data Table Example
#
# Update the Employee sheet's Username/Password sample data and mark the
# Password column (which already carries mailto hyperlinks) with the
# built-in "Hyperlink" cell style, then leave the selection where the
# author last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee")

# Expand the sample "Password" values
$ws.Range("E2").Value = "Secret@123!!ABC"
$ws.Range("E3").Value = "Secret@123!!ABC"

# Expand the sample "Username" values
$ws.Range("D2").Value = "JohnSmith2345005"
$ws.Range("D3").Value = "MaryAnn3845500"

# The Password cells are hyperlinked - give them the standard Hyperlink look
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"

# Leave the cursor where the author finished up
$ws.Range("D9").Select() | Out-Null
